$d = $word.ActiveDocument

# Find the "Author" style paragraph that contains "Edison Achalma" (the byline,
# right after the title heading), and insert a new "Author" style paragraph
# right after it containing the institutional affiliation line.
$target = $null
foreach ($p in $d.Paragraphs) {
    $style = $p.Range.ParagraphStyle
    $text = $p.Range.Text
    if ($style.NameLocal -eq "Author" -and $text.Trim() -eq "Edison Achalma") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range.Duplicate
    $r.Collapse(0)  # wdCollapseEnd - position right after the paragraph mark
    $r.InsertAfter("Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga`r")
    $newPara = $target.Next()
    $newPara.Style = "Author"
}
